$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 55285.49
$ws.Range("J17").Value = 56724.58
$ws.Range("L17").Value = 170173.74
$ws.Range("N17").Value = -170509.74
$ws.Range("H64").Value = 2990.75
$ws.Range("I64").Value = 3031.125
$ws.Range("J64").Value = 2910
$ws.Range("K64").Value = 3031.125
$ws.Range("L64").Value = 2910
$ws.Range("M64").Value = -2783.125
$ws.Range("N64").Value = -3406
$ws.Range("H67").Value = 2990.75
$ws.Range("I67").Value = 3031.125
$ws.Range("J67").Value = 2910
$ws.Range("K67").Value = 3031.125
$ws.Range("L67").Value = 2910
$ws.Range("M67").Value = -2173.125
$ws.Range("N67").Value = -4626
$ws.Range("H74").Value = 4750
$ws.Range("J74").Value = 5000
$ws.Range("L74").Value = 5000
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 4750
$ws.Range("J77").Value = 5000
$ws.Range("L77").Value = 25000
$ws.Range("N77").Value = -34360
$ws.Range("H141").Value = 465478.47
$ws.Range("I141").Value = 1998.2858
$ws.Range("J141").Value = 760420.4399999999
$ws.Range("K141").Value = 5994.857400000001
$ws.Range("L141").Value = 2281261.32
$ws.Range("M141").Value = -814.8574000000008
$ws.Range("N141").Value = -2291621.32

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1195.4286
$ws.Range("I74").Value = 1226.8334
$ws.Range("K74").Value = 1226.8334
$ws.Range("M74").Value = -352.8334
$ws.Range("H77").Value = 1195.4286
$ws.Range("I77").Value = 1226.8334
$ws.Range("K77").Value = 6134.166999999999
$ws.Range("M77").Value = -1766.166999999999
$ws.Range("H102").Value = 2957.6924
$ws.Range("I102").Value = 2313.6365
$ws.Range("K102").Value = 2313.6365
$ws.Range("M102").Value = -691.6365000000001
$ws.Range("H132").Value = 2682.5715
$ws.Range("I132").Value = 2550.1143
$ws.Range("J132").Value = 3344.8572
$ws.Range("K132").Value = 7650.342900000001
$ws.Range("L132").Value = 10034.5716
$ws.Range("M132").Value = -5120.342900000001
$ws.Range("N132").Value = -15094.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1756.3158
$ws.Range("I105").Value = 1657.0588
$ws.Range("K105").Value = 1657.0588
$ws.Range("M105").Value = 89.94119999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1726162
$ws.Range("I31").Value = 2084717.8
$ws.Range("J31").Value = 5094.9
$ws.Range("K31").Value = 2084717.8
$ws.Range("L31").Value = 5094.9
$ws.Range("M31").Value = -2084422.8
$ws.Range("N31").Value = -5684.9
$ws.Range("H34").Value = 1726162
$ws.Range("I34").Value = 2084717.8
$ws.Range("J34").Value = 5094.9
$ws.Range("K34").Value = 2084717.8
$ws.Range("L34").Value = 5094.9
$ws.Range("M34").Value = -2084515.8
$ws.Range("N34").Value = -5498.9
$ws.Range("H99").Value = 2444.1
$ws.Range("I99").Value = 1712.1333
$ws.Range("J99").Value = 4640
$ws.Range("K99").Value = 1712.1333
$ws.Range("L99").Value = 4640
$ws.Range("M99").Value = -214.1333
$ws.Range("N99").Value = -7636
$ws.Range("H122").Value = 2132.0527
$ws.Range("I122").Value = 1763.9333
$ws.Range("K122").Value = 5291.7999
$ws.Range("M122").Value = -2841.7999
$ws.Range("H126").Value = 2444.1
$ws.Range("I126").Value = 1712.1333
$ws.Range("J126").Value = 4640
$ws.Range("K126").Value = 5136.3999
$ws.Range("L126").Value = 13920
$ws.Range("M126").Value = -2666.3999
$ws.Range("N126").Value = -18860
$ws.Range("H133").Value = 8947.143
$ws.Range("J133").Value = 8947.143
$ws.Range("L133").Value = 8947.143
$ws.Range("N133").Value = -14007.143
$ws.Range("H134").Value = 1911.091
$ws.Range("I134").Value = 707.8823
$ws.Range("J134").Value = 6002
$ws.Range("K134").Value = 2123.6469
$ws.Range("L134").Value = 18006
$ws.Range("M134").Value = 411.3531000000003
$ws.Range("N134").Value = -23076

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 2602.8572
$ws.Range("I102").Value = 300
$ws.Range("J102").Value = 2986.6667
$ws.Range("K102").Value = 900
$ws.Range("L102").Value = 8960.000100000001
$ws.Range("M102").Value = 1534
$ws.Range("N102").Value = -13828.0001
$ws.Range("H103").Value = 2579.5806
$ws.Range("I103").Value = 782.6667
$ws.Range("J103").Value = 2772.1072
$ws.Range("K103").Value = 2348.0001
$ws.Range("L103").Value = 8316.321599999999
$ws.Range("M103").Value = -1469.0001
$ws.Range("N103").Value = -10074.3216

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2195.353
$ws.Range("I113").Value = 1665.5454
$ws.Range("J113").Value = 3166.6667
$ws.Range("K113").Value = 1665.5454
$ws.Range("L113").Value = 3166.6667
$ws.Range("M113").Value = 504.4546
$ws.Range("N113").Value = -7506.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6252581.5
$ws.Range("I7").Value = 20002798
$ws.Range("J7").Value = 2483
$ws.Range("K7").Value = 20002798
$ws.Range("L7").Value = 2483
$ws.Range("M7").Value = -20002686
$ws.Range("N7").Value = -2707
$ws.Range("H40").Value = 6283.4614
$ws.Range("I40").Value = 7298.3335
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 7298.3335
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -7162.3335
$ws.Range("N40").Value = -4272
$ws.Range("H122").Value = 3143.7778
$ws.Range("I122").Value = 2371.6365
$ws.Range("J122").Value = 4357.143
$ws.Range("K122").Value = 7114.9095
$ws.Range("L122").Value = 13071.429
$ws.Range("M122").Value = -4664.9095
$ws.Range("N122").Value = -17971.429
$ws.Range("H126").Value = 6252581.5
$ws.Range("I126").Value = 20002798
$ws.Range("J126").Value = 2483
$ws.Range("K126").Value = 60008394
$ws.Range("L126").Value = 7449
$ws.Range("M126").Value = -60005924
$ws.Range("N126").Value = -12389
$ws.Range("H132").Value = 2821.92
$ws.Range("I132").Value = 3049.6
$ws.Range("J132").Value = 2765
$ws.Range("K132").Value = 9148.799999999999
$ws.Range("L132").Value = 8295
$ws.Range("M132").Value = -6618.799999999999
$ws.Range("N132").Value = -13355

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1993.8572
$ws.Range("I96").Value = 1767.6666
$ws.Range("J96").Value = 2163.5
$ws.Range("K96").Value = 1767.6666
$ws.Range("L96").Value = 2163.5
$ws.Range("M96").Value = -394.6666
$ws.Range("N96").Value = -4909.5
$ws.Range("H122").Value = 271736.75
$ws.Range("I122").Value = 313526.75
$ws.Range("J122").Value = 4280.8
$ws.Range("K122").Value = 940580.25
$ws.Range("L122").Value = 12842.4
$ws.Range("M122").Value = -938130.25
$ws.Range("N122").Value = -17742.4
$ws.Range("H126").Value = 3848396
$ws.Range("I126").Value = 1493.0625
$ws.Range("J126").Value = 10003441
$ws.Range("K126").Value = 4479.1875
$ws.Range("L126").Value = 30010323
$ws.Range("M126").Value = -2009.1875
$ws.Range("N126").Value = -30015263
